$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (TB), C (d2S), D (K), E (IP) and G (sum) for rows 2-14.
# F (Win) and A (date) are unchanged. G = B + C + D + E (Win excluded from sum).
$data = @{
    2  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    3  = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548)
    4  = @(0.1190320826869504, 0.04071648406533734, 0.7527432677738641, 0.4942365360607697, 1.406728370586922)
    5  = @(1.455362044514542, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 3.009163075608874)
    6  = @(0.1190320826869504, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 2.418468675747795)
    7  = @(0.6606524410359556, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 2.960089034096801)
    8  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    9  = @(0.6606524410359556, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 6.348428708163715)
    10 = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 4.358119930609447)
    11 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    12 = @(0.01293466051926884, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 2.915692546614173)
    13 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548)
    14 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 10.19245300693656, 18.67282528286833)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 7).Value = $vals[4]
}
